# Generate Report for Handback
# Update timestamps / priority values that are shared across the
# Overview / zh-cn / de-de sheets. Because several cells point at the
# *same* shared string, every cell sharing an old value must be rewritten
# to the new value so the workbook's shared-string table stays correctly
# deduplicated after save (otherwise only the first-touched cell would
# pick up a brand-new string while its siblings kept the stale text).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), both data rows
# shared the same timestamp "2016-09-06 16:10:29" -> "2016-09-06 16:13:25"
$wsOverview.Range("G2").Value = "2016-09-06 16:13:25"
$wsOverview.Range("G3").Value = "2016-09-06 16:13:25"

# zh-cn sheet
#  Priority column (E): "ht" -> "mt" for both rows
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

#  Correspond Handoff Datetime column (H): "2016-09-06 16:09:57" -> "2016-09-06 16:13:08"
$wsZhCn.Range("H2").Value = "2016-09-06 16:13:08"
$wsZhCn.Range("H3").Value = "2016-09-06 16:13:08"

#  Correspond Handback DateTime column (K): "2016-09-06 16:11:45" -> "2016-09-06 16:14:14"
$wsZhCn.Range("K2").Value = "2016-09-06 16:14:14"
$wsZhCn.Range("K3").Value = "2016-09-06 16:14:14"

# de-de sheet
#  Priority column (E): "ht" -> "mt" for both rows
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

#  Correspond Handoff Datetime column (H) shares the Overview timestamp:
#  "2016-09-06 16:10:29" -> "2016-09-06 16:13:25"
$wsDeDe.Range("H2").Value = "2016-09-06 16:13:25"
$wsDeDe.Range("H3").Value = "2016-09-06 16:13:25"

#  Correspond Handback DateTime column (K): "2016-09-06 16:12:10" -> "2016-09-06 16:14:33"
$wsDeDe.Range("K2").Value = "2016-09-06 16:14:33"
$wsDeDe.Range("K3").Value = "2016-09-06 16:14:33"
